$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing "BAYMAG - Coretop" header from V1 to X1, and add new headers in V1/W1
$ws.Range("X1").Value = $ws.Range("V1").Value()
$ws.Range("V1").Value = "MgCa Coretop modelled temperature"
$ws.Range("W1").Value = "MgCa Temperature anomaly_Original - Coretop"

# Copy the header style (bold/centered) from an existing header cell to the new header cells
$ws.Range("U1").Copy()
$ws.Range("V1:X1").PasteSpecial(-4122)  # xlPasteFormats

# Update row 2 values
$ws.Range("S2").Value = 22.22
$ws.Range("T2").Value = -1.791507281679124
$ws.Range("U2").Value = 2.938683627411777
$ws.Range("X2").Value = 0.6081636400000008
$ws.Range("V2").Value = 24.5493
$ws.Range("W2").Value = -4.12202727
